# Release Log Form - F6.xlsx : "feat: sops Update 4"
#
# Moves this form from the "Software Service Catalog" area to the
# "Software Development Lifecycle" area: renames the visible sheet,
# repoints the Print_Area defined name at the new sheet name, un-hides
# the helper sheet, nudges the saved scroll position, and updates the
# revision date stamped in the footer.

$wb = $excel.ActiveWorkbook

# --- Sheets -----------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the main sheet S-SW-SC-06 -> F-SW-SD-06
$ws1.Name = "F-SW-SD-06"

# Un-hide "Sheet2" (was state="hidden")
$ws2.Visible = $true

# --- Defined name / print area ----------------------------------------
# The Print_Area name tracks the renamed sheet by key automatically, but
# its RefersTo text still points at the old sheet title - fix it up.
$printArea = $wb.Names.Item("F-SW-SD-06!Print_Area")
$printArea.RefersTo = "='F-SW-SD-06'!`$A`$1:`$G`$31"

# --- View state ---------------------------------------------------------
$ws1.Activate()
[void]$ws1.Range("A6").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1

# Restore the originally-selected cell on the sheet
[void]$ws1.Range("E24").Select()

# --- Footer -------------------------------------------------------------
# Rev date placeholder (0/0/2025) -> actual revision date (01/10/2025)
$ps = $ws1.PageSetup
$ps.LeftFooter = "&14Issue No.: (1)"
$ps.CenterFooter = "&14F-SW-SD/06"
$ps.RightFooter = "&14Rev:0(01/10/2025)"
